$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shelby's updates: fill in the remaining Functional Components cells for
# rows 4 and 5 (the "Additional motorcycle parking" user stories).
$ws.Range("H4").Value = "1. Pop up options.`n2. Voice interaction."
$ws.Range("G5").Value = "1. Additional routes shown that the user didn’t select.`n2.Show impediments on any route."
$ws.Range("H5").Value = "1.A small window.`n2. A free text area to accept user input for this field."

# Restore the author's view state (scroll position / selection) after the edits.
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 1
$ws.Range("C7").Select()
